$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers / values, entered in the same order the author typed them ---
# Column D: CreateAccount
$ws.Range("D1").Value = "CreateAccount"
$ws.Range("D2").Value = "Account  acc = New Account(Name = 'Account Automation 1', Type = 'Analyst');insert acc ;"

# Column E: fetchaccount
$ws.Range("E1").Value = "fetchaccount"
$ws.Range("E2").Value = "Select Name , Id from Account where Name = 'Account Automation 1' AND Createdby.Id = '005q0000003GGfP' Order by CreatedDate DESC Limit 1"

# Column F: CreateServiceContract
$ws.Range("F1").Value = "CreateServiceContract"
# Column G: FetchServiceContract
$ws.Range("G1").Value = "FetchServiceContract"

$ws.Range("F2").Value = "SVMXC__Service_Contract__c service_contract = New SVMXC__Service_Contract__c (Name = 'SCON RS_1022', SVMXC__Active__c = true , SVMXC__All_Contacts_Covered__c = true , SVMXC__Company__c = '001q000000kxZfw');insert service_contract;"
$ws.Range("G2").Value = "Select Name , Id from SVMXC__Service_Contract__c where Name = 'SCON RS_1022' AND Createdby.Id = '005q0000003GGfP' Order by CreatedDate DESC Limit 1"

$ws.Range("H2").Value = "SVMXC"
$ws.Range("I1").Value = "Account"
$ws.Range("I2").Value = "001q000000kxZfw"
$ws.Range("J1").Value = "Username"
# J2 keeps a leading apostrophe in the original workbook (quote-prefixed text)
$ws.Range("J2").Value = "'005q0000003GGfP"
$ws.Range("H1").Value = "OrgDetails"

# --- Formatting: the "Select ..." query cells (E2, G2) reuse the Monaco/blue
# query style already used by A2/B2/C2; copy that formatting across. ---
$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths for the new columns D:J ---
$ws.Columns.Item(4).ColumnWidth = 65.166666667
$ws.Columns.Item(5).ColumnWidth = 44
$ws.Columns.Item(6).ColumnWidth = 77.5
$ws.Columns.Item(7).ColumnWidth = 80.666666667
$ws.Columns.Item(8).ColumnWidth = 35.5
$ws.Columns.Item(9).ColumnWidth = 37
$ws.Columns.Item(10).ColumnWidth = 27.333333333

# --- View: scroll so column F is left-most and select G2 (matches the
# author's on-screen state when the change was saved) ---
$ws.Range("G2").Select()
$excel.ActiveWindow.ScrollColumn = 6
